{"js": "// Update the date line and the 25 division-fact answers in the table.\n// The document has one intro paragraph with the date, followed by a single\n// 5-column table whose populated rows are 0, 4, 8, 12, 16 (the rows in\n// between are blank spacer rows).\n\nconst body = context.document.body;\n\n// --- 1. Update the date paragraph -----------------------------------\nconst paras = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\nconst datePara = paras.items[0];\ndatePara.load(\"text\");\nawait context.sync();\nif (datePara.text.trim() === \"2025-12-13 Saturday\") {\n  datePara.insertText(\"2025-12-14 Sunday\", \"Replace\");\n}\n\n// --- 2. Update the table cell contents --------------------------------\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Row (table-relative, 0-based) -> column -> new text.\nconst rowUpdates = {\n  0: [\"52\u00f76=8, 4\", \"79\u00f78=9, 7\", \"11\u00f79=1, 2\", \"95\u00f79=10, 5\", \"82\u00f76=13, 4\"],\n  4: [\"47\u00f78=5, 7\", \"78\u00f73=26, 0\", \"50\u00f78=6, 2\", \"78\u00f73=26, 0\", \"20\u00f73=6, 2\"],\n  8: [\"90\u00f73=30, 0\", \"21\u00f72=10, 1\", \"94\u00f73=31, 1\", \"17\u00f75=3, 2\", \"14\u00f74=3, 2\"],\n  12: [\"92\u00f72=46, 0\", \"60\u00f73=20, 0\", \"89\u00f74=22, 1\", \"57\u00f79=6, 3\", \"27\u00f78=3, 3\"],\n  16: [\"17\u00f76=2, 5\", \"18\u00f78=2, 2\", \"44\u00f75=8, 4\", \"42\u00f78=5, 2\", \"51\u00f79=5, 6\"]\n};\n\nfor (const rowIdx of Object.keys(rowUpdates)) {\n  const r = Number(rowIdx);\n  const values = rowUpdates[rowIdx];\n  for (let c = 0; c < values.length; c++) {\n    table.getCell(r, c).value = values[c];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the 25 division-fact answers in the table.\n# The document has one intro paragraph with the date, followed by a single\n# 5-column table whose populated rows are COM rows 1, 5, 9, 13, 17 (the\n# rows in between are blank spacer rows). Table.Cell(row, col) is 1-based.\n\n$d = $word.ActiveDocument\n\n# --- 1. Update the date paragraph -----------------------------------\n$datePara = $d.Paragraphs(1)\nif ($datePara.Range.Text.Trim() -eq \"2025-12-13 Saturday\") {\n    $datePara.Range.Text = \"2025-12-14 Sunday\"\n}\n\n# --- 2. Update the table cell contents --------------------------------\n$t = $d.Tables(1)\n\n$rowValues = @{\n    1  = @(\"52\u00f76=8, 4\", \"79\u00f78=9, 7\", \"11\u00f79=1, 2\", \"95\u00f79=10, 5\", \"82\u00f76=13, 4\")\n    5  = @(\"47\u00f78=5, 7\", \"78\u00f73=26, 0\", \"50\u00f78=6, 2\", \"78\u00f73=26, 0\", \"20\u00f73=6, 2\")\n    9  = @(\"90\u00f73=30, 0\", \"21\u00f72=10, 1\", \"94\u00f73=31, 1\", \"17\u00f75=3, 2\", \"14\u00f74=3, 2\")\n    13 = @(\"92\u00f72=46, 0\", \"60\u00f73=20, 0\", \"89\u00f74=22, 1\", \"57\u00f79=6, 3\", \"27\u00f78=3, 3\")\n    17 = @(\"17\u00f76=2, 5\", \"18\u00f78=2, 2\", \"44\u00f75=8, 4\", \"42\u00f78=5, 2\", \"51\u00f79=5, 6\")\n}\n\nforeach ($row in $rowValues.Keys) {\n    $values = $rowValues[$row]\n    for ($c = 0; $c -lt $values.Length; $c++) {\n        $t.Cell($row, $c + 1).Range.Text = $values[$c]\n    }\n}\n"}
